$d = $word.ActiveDocument

# "Core Competencies" table is the first table in the resume.
$t = $d.Tables.Item(1)

# -- Row 7, Col 1: "sass" -> "MVC" -----------------------------------------
$sassCell  = $t.Cell(7, 1)
$sassRange = $d.Range($sassCell.Range.Start, $sassCell.Range.Start + 4)
$sassRange.Text = "MVC"

# -- Row 7, Col 3: remove "Vue.js", leaving a bare empty paragraph ---------
$vueCell      = $t.Cell(7, 3)
$vueTextRange = $d.Range($vueCell.Range.Start, $vueCell.Range.Start + 6)
$vueTextRange.Delete()
$vueCell.Range.Paragraphs.Item(1).Style = "Normal"
